# Insert two new weekly price rows for "Betarraga" (Macroferia Regional de
# Talca) ahead of the existing row 480, shifting the old rows 480-514 down
# to 482-516.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows starting at row 480 (pushes existing data down).
$ws.Range("A480:R481").EntireRow.Insert()

# New row 480
$ws.Range("A480").Value = 5
$ws.Range("B480").Value = "Macroferia Regional de Talca"
$ws.Range("C480").Value = "Maule"
$ws.Range("D480").Value = 45021
$ws.Range("D480").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E480").Value = 7
$ws.Range("F480").Value = 100114014
$ws.Range("G480").Value = "Betarraga"
$ws.Range("H480").Value = "Sin especificar"
$ws.Range("I480").Value = "Primera"
$ws.Range("J480").Value = 3000
$ws.Range("K480").Value = 600
$ws.Range("L480").Value = 600
$ws.Range("M480").Value = 600
$ws.Range("N480").Value = "`$/paquete 5 unidades"
$ws.Range("O480").Value = "Región del Maule"
$ws.Range("P480").Value = 120
$ws.Range("Q480").Value = 5
$ws.Range("R480").Value = "Hortaliza"

# New row 481
$ws.Range("A481").Value = 5
$ws.Range("B481").Value = "Macroferia Regional de Talca"
$ws.Range("C481").Value = "Maule"
$ws.Range("D481").Value = 45021
$ws.Range("D481").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E481").Value = 7
$ws.Range("F481").Value = 100114014
$ws.Range("G481").Value = "Betarraga"
$ws.Range("H481").Value = "Sin especificar"
$ws.Range("I481").Value = "Segunda"
$ws.Range("J481").Value = 3000
$ws.Range("K481").Value = 500
$ws.Range("L481").Value = 500
$ws.Range("M481").Value = 500
$ws.Range("N481").Value = "`$/paquete 5 unidades"
$ws.Range("O481").Value = "Región del Maule"
$ws.Range("P481").Value = 100
$ws.Range("Q481").Value = 5
$ws.Range("R481").Value = "Hortaliza"
